$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.00863423648946369, 1.081231980802479, 8.112271609790723, 2.848204980297367, 2.87966458212119, 46),
    @(0.1079425770541491, 1.025681875173506, 5.546319209064285, 2.355062463941092, 2.379171181174347, 45),
    @(0.08673112263663806, 0.8989375388443066, 4.47623325838249, 2.115711052668225, 2.138371923145345, 44),
    @(0.1560044769001684, 0.9224880001651093, 4.879102887750935, 2.20886914228773, 2.229429331798652, 43),
    @(0.1600298176880285, 0.9424896015931804, 4.916621931739378, 2.217345695136277, 2.238371143106592, 42),
    @(0.1728584914046674, 0.9206488348686618, 4.635874527679912, 2.153108108683796, 2.172819442882898, 41),
    @(0.1929609596835913, 0.9736223895218249, 4.804027380007986, 2.191809156840072, 2.211112598008703, 40),
    @(0.1960041159357992, 0.9531248790116683, 4.870747427708349, 2.206976988486366, 2.226992678098346, 39),
    @(0.2155726075138289, 0.9906079526445191, 5.010930836497054, 2.238510852441206, 2.258015469991629, 38),
    @(0.2079488633544275, 0.9496862135200643, 5.054960332149459, 2.248323893959555, 2.269566474070086, 37)
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 2
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}

$wb.Save()
